# RoughPad.xlsx: add a new "Sheet2" after "Sheet1", populate it with the
# Agent record field labels, and make it the active/selected sheet.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1 (Add() defaults to "before active
# sheet", so pass Sheet1 explicitly as the After-anchor).
$sheet2 = $wb.Worksheets.Add($null, $sheet1)

# Fill column A with the field-label header rows.
$labels = @("Agent", "Name", "Gender", "Birthday", "Mood", "Family")
for ($i = 0; $i -lt $labels.Length; $i++) {
    $sheet2.Cells.Item($i + 1, 1).Value = $labels[$i]
}

# Match the committed selection/active-tab state: Sheet2 is selected, with
# A2 the active cell.
$sheet2.Range("A2").Select()
